$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 49.07461933333334
$ws.Range("H2").Value = 147.223858
$ws.Range("I2").Value = 0.2082521693470354
$ws.Range("J2").Value = 0.2082521693470354
$ws.Range("M2").Value = 3.317411333333334
$ws.Range("N2").Value = 9.952234000000001
$ws.Range("O2").Value = 0.1227370957132105
$ws.Range("P2").Value = 0.1227370957132105
$ws.Range("Q2").Value = 162.8006983554191
$ws.Range("R2").Value = 1465.206285198772
$ws.Range("S2").Value = 0.02556026644163081
$ws.Range("T2").Value = 0.02556026644163081
$ws.Range("G3").Value = 49.07461933333334
$ws.Range("H3").Value = 147.223858
$ws.Range("I3").Value = 0.2082521693470354
$ws.Range("J3").Value = 0.2082521693470354
$ws.Range("O3").Value = 0.4742843403616469
$ws.Range("P3").Value = 0.4742843403616469
$ws.Range("Q3").Value = 629.0993067844336
$ws.Range("R3").Value = 5661.893761059902
$ws.Range("S3").Value = 0.09877074276764067
$ws.Range("T3").Value = 0.09877074276764064
$ws.Range("G4").Value = 49.07461933333334
$ws.Range("H4").Value = 147.223858
$ws.Range("I4").Value = 0.2082521693470354
$ws.Range("J4").Value = 0.2082521693470354
$ws.Range("M4").Value = 7.649101666666667
$ws.Range("N4").Value = 22.947305
$ws.Range("O4").Value = 0.2830003364214742
$ws.Range("P4").Value = 0.2830003364214742
$ws.Range("Q4").Value = 375.3767525336323
$ws.Range("R4").Value = 3378.39077280269
$ws.Range("S4").Value = 0.05893543398571284
$ws.Range("T4").Value = 0.05893543398571283
$ws.Range("G5").Value = 49.07461933333334
$ws.Range("H5").Value = 147.223858
$ws.Range("I5").Value = 0.2082521693470354
$ws.Range("J5").Value = 0.2082521693470354
$ws.Range("M5").Value = 3.242843
$ws.Range("N5").Value = 9.728529
$ws.Range("O5").Value = 0.1199782275036684
$ws.Range("P5").Value = 0.1199782275036684
$ws.Range("Q5").Value = 159.1412857827647
$ws.Range("R5").Value = 1432.271572044882
$ws.Range("S5").Value = 0.0249857261520511
$ws.Range("T5").Value = 0.02498572615205109
$ws.Range("I6").Value = 0.4723739808466805
$ws.Range("J6").Value = 0.4723739808466804
$ws.Range("M6").Value = 3.317411333333334
$ws.Range("N6").Value = 9.952234000000001
$ws.Range("O6").Value = 0.1227370957132105
$ws.Range("P6").Value = 0.1227370957132105
$ws.Range("Q6").Value = 369.2773727538782
$ws.Range("R6").Value = 3323.496354784904
$ws.Range("S6").Value = 0.05797781049960929
$ws.Range("T6").Value = 0.05797781049960929
$ws.Range("I7").Value = 0.4723739808466805
$ws.Range("J7").Value = 0.4723739808466804
$ws.Range("O7").Value = 0.4742843403616469
$ws.Range("P7").Value = 0.4742843403616469
$ws.Range("S7").Value = 0.2240395819098731
$ws.Range("T7").Value = 0.224039581909873
$ws.Range("I8").Value = 0.4723739808466805
$ws.Range("J8").Value = 0.4723739808466804
$ws.Range("M8").Value = 7.649101666666667
$ws.Range("N8").Value = 22.947305
$ws.Range("O8").Value = 0.2830003364214742
$ws.Range("P8").Value = 0.2830003364214742
$ws.Range("Q8").Value = 851.4591298980644
$ws.Range("R8").Value = 7663.13216908258
$ws.Range("S8").Value = 0.1336819954963616
$ws.Range("T8").Value = 0.1336819954963616
$ws.Range("I9").Value = 0.4723739808466805
$ws.Range("J9").Value = 0.4723739808466804
$ws.Range("M9").Value = 3.242843
$ws.Range("N9").Value = 9.728529
$ws.Range("O9").Value = 0.1199782275036684
$ws.Range("P9").Value = 0.1199782275036684
$ws.Range("Q9").Value = 360.9768047937693
$ws.Range("R9").Value = 3248.791243143924
$ws.Range("S9").Value = 0.05667459294083654
$ws.Range("T9").Value = 0.05667459294083654
$ws.Range("G10").Value = 40.34671033333333
$ws.Range("H10").Value = 121.040131
$ws.Range("I10").Value = 0.1712145721571795
$ws.Range("J10").Value = 0.1712145721571795
$ws.Range("M10").Value = 3.317411333333334
$ws.Range("N10").Value = 9.952234000000001
$ws.Range("O10").Value = 0.1227370957132105
$ws.Range("P10").Value = 0.1227370957132105
$ws.Range("Q10").Value = 133.8466341225171
$ws.Range("R10").Value = 1204.619707102654
$ws.Range("S10").Value = 0.02101437933035213
$ws.Range("T10").Value = 0.02101437933035213
$ws.Range("G11").Value = 40.34671033333333
$ws.Range("H11").Value = 121.040131
$ws.Range("I11").Value = 0.1712145721571795
$ws.Range("J11").Value = 0.1712145721571795
$ws.Range("O11").Value = 0.4742843403616469
$ws.Range("P11").Value = 0.4742843403616469
$ws.Range("Q11").Value = 517.2141495245766
$ws.Range("R11").Value = 4654.927345721188
$ws.Range("S11").Value = 0.08120439041586947
$ws.Range("T11").Value = 0.08120439041586945
$ws.Range("G12").Value = 40.34671033333333
$ws.Range("H12").Value = 121.040131
$ws.Range("I12").Value = 0.1712145721571795
$ws.Range("J12").Value = 0.1712145721571795
$ws.Range("M12").Value = 7.649101666666667
$ws.Range("N12").Value = 22.947305
$ws.Range("O12").Value = 0.2830003364214742
$ws.Range("P12").Value = 0.2830003364214742
$ws.Range("Q12").Value = 308.6160892552172
$ws.Range("R12").Value = 2777.544803296955
$ws.Range("S12").Value = 0.04845378152074057
$ws.Range("T12").Value = 0.04845378152074056
$ws.Range("G13").Value = 40.34671033333333
$ws.Range("H13").Value = 121.040131
$ws.Range("I13").Value = 0.1712145721571795
$ws.Range("J13").Value = 0.1712145721571795
$ws.Range("M13").Value = 3.242843
$ws.Range("N13").Value = 9.728529
$ws.Range("O13").Value = 0.1199782275036684
$ws.Range("P13").Value = 0.1199782275036684
$ws.Range("Q13").Value = 130.8380471774777
$ws.Range("R13").Value = 1177.542424597299
$ws.Range("S13").Value = 0.02054202089021734
$ws.Range("T13").Value = 0.02054202089021733
$ws.Range("G14").Value = 34.91373066666667
$ws.Range("H14").Value = 104.741192
$ws.Range("I14").Value = 0.1481592776491046
$ws.Range("J14").Value = 0.1481592776491046
$ws.Range("M14").Value = 3.317411333333334
$ws.Range("N14").Value = 9.952234000000001
$ws.Range("O14").Value = 0.1227370957132105
$ws.Range("P14").Value = 0.1227370957132105
$ws.Range("Q14").Value = 115.8232058025476
$ws.Range("R14").Value = 1042.408852222928
$ws.Range("S14").Value = 0.01818463944161828
$ws.Range("T14").Value = 0.01818463944161828
$ws.Range("G15").Value = 34.91373066666667
$ws.Range("H15").Value = 104.741192
$ws.Range("I15").Value = 0.1481592776491046
$ws.Range("J15").Value = 0.1481592776491046
$ws.Range("O15").Value = 0.4742843403616469
$ws.Range("P15").Value = 0.4742843403616469
$ws.Range("Q15").Value = 447.5674810734499
$ws.Range("R15").Value = 4028.107329661048
$ws.Range("S15").Value = 0.07026962526826366
$ws.Range("T15").Value = 0.07026962526826365
$ws.Range("G16").Value = 34.91373066666667
$ws.Range("H16").Value = 104.741192
$ws.Range("I16").Value = 0.1481592776491046
$ws.Range("J16").Value = 0.1481592776491046
$ws.Range("M16").Value = 7.649101666666667
$ws.Range("N16").Value = 22.947305
$ws.Range("O16").Value = 0.2830003364214742
$ws.Range("P16").Value = 0.2830003364214742
$ws.Range("Q16").Value = 267.0586754319511
$ws.Range("R16").Value = 2403.52807888756
$ws.Range("S16").Value = 0.04192912541865921
$ws.Range("T16").Value = 0.0419291254186592
$ws.Range("G17").Value = 34.91373066666667
$ws.Range("H17").Value = 104.741192
$ws.Range("I17").Value = 0.1481592776491046
$ws.Range("J17").Value = 0.1481592776491046
$ws.Range("M17").Value = 3.242843
$ws.Range("N17").Value = 9.728529
$ws.Range("O17").Value = 0.1199782275036684
$ws.Range("P17").Value = 0.1199782275036684
$ws.Range("Q17").Value = 113.2197470962854
$ws.Range("R17").Value = 1018.977723866568
$ws.Range("S17").Value = 0.01777588752056345
$ws.Range("T17").Value = 0.01777588752056345
